$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.238.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.703.01"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.83%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.698.27"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.89%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("E10").Value = "  +5.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.03%  "

$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.13"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.322.18"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.707.62"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.147.61"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.20"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  -1.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.68"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.33%  "

$ws.Range("E25").Value = "  +3.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.91%  "

$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.42"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.845.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.647.58"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.91%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.77"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("E40").Value = "  -3.33%  "

$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.10"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.51"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("E44").Value = "  -1.80%  "

$ws.Range("E45").Value = "  -1.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.34"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.757.29"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.10%  "

$ws.Range("E51").Value = "  -0.75%  "

